$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: the "b.md" row now reports "Ready for handoff" for both
# the zh-cn and de-de columns (previously "Handed back: in sync with en-US").
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 (b.md) now has a fresh handoff file/date, and status
# changes to "Ready for handoff".
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "2016-02-22 17:32:26"

# The "Latest Handoff File" cell (C3) needs its displayed text updated to
# the new handoff file name while keeping the same underlying hyperlink
# target. Rebuild every hyperlink on the sheet (in original order) so the
# relationship ids line up the same way, since this host cannot mutate a
# single hyperlink's display text in place.
$zhLinks = @(
    @{ Cell = "A2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/87aed72acb7ccea224e9f7e2ce18451bee58e873/e2e/a.md"; Display = "a.md" },
    @{ Cell = "C2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/90df88de8e333fc883e877a686133e687b711794/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Cell = "E2"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/88358cd83a62ffe9ddf5c12d95f2268525f7e092/e2e/a.md"; Display = "a.md" },
    @{ Cell = "F2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8fdabbf213d49c76ef4d18867df1ff6aef0a7009/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Cell = "A3"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/87aed72acb7ccea224e9f7e2ce18451bee58e873/e2e/b.md"; Display = "b.md" },
    @{ Cell = "C3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/90df88de8e333fc883e877a686133e687b711794/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf" },
    @{ Cell = "E3"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/88358cd83a62ffe9ddf5c12d95f2268525f7e092/e2e/a.md"; Display = "a.md" },
    @{ Cell = "F3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8fdabbf213d49c76ef4d18867df1ff6aef0a7009/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Cell = "A4"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/87aed72acb7ccea224e9f7e2ce18451bee58e873/.localization-config"; Display = ".localization-config" }
)

$wsZh.Hyperlinks.Delete()
foreach ($link in $zhLinks) {
    $wsZh.Hyperlinks.Add($wsZh.Range($link.Cell), $link.Address, "", "", $link.Display) | Out-Null
}

# ---------------------------------------------------------------------------
# de-de sheet: same kind of update as zh-cn, with the de-de handoff file.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "2016-02-22 17:32:37"

$deLinks = @(
    @{ Cell = "A2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/87aed72acb7ccea224e9f7e2ce18451bee58e873/e2e/a.md"; Display = "a.md" },
    @{ Cell = "C2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/efec60cb81a854ced46235ad54569a8ab6fbc0fd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Cell = "E2"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/74948b4db71264ca0d434087994ab80c76cf48c6/e2e/a.md"; Display = "a.md" },
    @{ Cell = "F2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/811b3825df1c2b482129a279efae658bf37f5487/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Cell = "A3"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/87aed72acb7ccea224e9f7e2ce18451bee58e873/e2e/b.md"; Display = "b.md" },
    @{ Cell = "C3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/efec60cb81a854ced46235ad54569a8ab6fbc0fd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf" },
    @{ Cell = "E3"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/74948b4db71264ca0d434087994ab80c76cf48c6/e2e/a.md"; Display = "a.md" },
    @{ Cell = "F3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/811b3825df1c2b482129a279efae658bf37f5487/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Cell = "A4"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/87aed72acb7ccea224e9f7e2ce18451bee58e873/.localization-config"; Display = ".localization-config" }
)

$wsDe.Hyperlinks.Delete()
foreach ($link in $deLinks) {
    $wsDe.Hyperlinks.Add($wsDe.Range($link.Cell), $link.Address, "", "", $link.Display) | Out-Null
}
